$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q1" sheet (position 2) to create the new
#    "2022-Q3" sheet, inserted *before* "2022-Q1" so the tab order becomes
#    总计, 2022-Q3, 2022-Q1, 2021-Q3. Duplicating (rather than Worksheets.Add)
#    keeps every header/column style identical to the sibling quarter sheet.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(2)
$q1.Copy($q1)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Give the new data row (row 3) the same look as row 2's index cell before
# filling in values.
$q3.Range("A2").Copy()
$q3.Range("A3").PasteSpecial(-4122)

# The fund-code / name / ratio columns in this workbook are stored as text
# (e.g. "006555" keeps its leading zero), so force text formatting before
# typing numeric-looking strings into B:G.
$q3.Range("B2:G3").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "006555"
$q3.Range("C2").Value = "浦银安盛全球智能科技股票（QDII）A"
$q3.Range("D2").Value = "0.25"
$q3.Range("E2").Value = "84.65"
$q3.Range("F2").Value = "2.16"
$q3.Range("G2").Value = "0.0054"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "014002"
$q3.Range("C3").Value = "浦银安盛全球智能科技股票（QDII）C"
$q3.Range("D3").Value = "0.01"
$q3.Range("E3").Value = "84.65"
$q3.Range("F3").Value = "2.16"
$q3.Range("G3").Value = "0.0002"
$q3.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert the new 2022-Q3 row at the
#    top of the data (row 2), push the existing 2022-Q1 row down to row 3,
#    and append a new row 4 with the 2021-Q3 totals.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Give row 4's index cell (col A) the same style as the existing rows above.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01

$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.07000000000000001

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.01

# ---------------------------------------------------------------------------
# 3) Restore the originally-active tab ("2021-Q3", the last sheet) — copying
#    a sheet shifts Excel's focus onto the newly created copy, so put the
#    selection back where the workbook had it before we started editing.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
